$d = $word.ActiveDocument

# 1) Kontakt section: "Antrimon" + " " runs merge into a single run "Antrimon "
#    (the proofErr spell-check bookmarks around "Antrimon" are removed as a
#    side effect of replacing across the run boundary)
$d.Content.Find.Execute("Antrimon Group AG", $true, $false, $false, $false, $false, $true, 1, $false, "Antrimon Group AG", 2)

# 2) "Surentalstrasse" + " 10, 6210 Sursee" runs merge into a single run
$d.Content.Find.Execute("Surentalstrasse 10, 6210 Sursee", $true, $false, $false, $false, $false, $true, 1, $false, "Surentalstrasse 10, 6210 Sursee", 2)

# 3) Table cell: "Antrimon" + " Group" runs merge into a single run
$d.Content.Find.Execute("Antrimon Group", $true, $false, $false, $false, $false, $true, 1, $false, "Antrimon Group", 2)

# 4) "CRMNr" inline text merges with surrounding runs into one run
$d.Content.Find.Execute("Zuerst in SuperOffice und anschliessend in SAP mit der Verkaufsnummer aus SuperOffice (CRMNr)?", $true, $false, $false, $false, $false, $true, 1, $false, "Zuerst in SuperOffice und anschliessend in SAP mit der Verkaufsnummer aus SuperOffice (CRMNr)?", 2)

# 5) The last (previously empty) paragraph after the Q&A gets "New version"
$qa = $d.Content
$qa.Find.Execute("Zuerst in SuperOffice und anschliessend in SAP mit der Verkaufsnummer aus SuperOffice (CRMNr)?")
$qaPara = $qa.Paragraphs(1)
$lastPara = $qaPara.Next().Next()
$lastPara.Range.Text = "New version"

# 6) Update cached SAVEDATE field text in each footer from 15.01.2023 to 16.01.2023
$section = $d.Sections(1)
for ($i = 1; $i -le $section.Footers.Count; $i++) {
    $footer = $section.Footers($i)
    if ($footer.Exists) {
        $footer.Range.Find.Execute("15.01.2023", $true, $false, $false, $false, $false, $true, 1, $false, "16.01.2023", 2)
    }
}
